$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    "B2" = 0.582312245481221
    "C2" = 2.31174574942257
    "D2" = 9.71846063957016
    "E2" = 8.44423096866041
    "F2" = 24.427811206171
    "G2" = 100.312929246159

    "B3" = 3.75402348695717
    "C3" = 12.3523395214045
    "D3" = 43.5372812714913
    "E3" = 38.4957567164209
    "F3" = 103.408058751747
    "G3" = 392.685850239176

    "B4" = 17.0056341062422
    "C4" = 59.1151022823119
    "D4" = 149.464993029533
    "E4" = 137.096776986925
    "F4" = 311.316952676775
    "G4" = 850.990988524279

    "B5" = 0.255033144765084
    "C5" = 0.932014984439093
    "D5" = 3.08947832460019
    "E5" = 2.69396590295945
    "F5" = 7.51138324483634
    "G5" = 25.0156050214649
}

foreach ($addr in $data.Keys) {
    $ws.Range($addr).Value = $data[$addr]
}
